$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.781.92"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "3.600.79"
$ws.Range("E3").Value = "  +1.79%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'603.65"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").Value = "'136.89"
$ws.Range("E6").Value = "  -1.51%  "
$ws.Range("D7").Value = "3.599.97"
$ws.Range("E7").Value = "  +1.76%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.498"
$ws.Range("E9").Value = "  +1.89%  "
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("D11").Value = "'7.25"
$ws.Range("E11").Value = "  +4.90%  "
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").Value = "4.213.52"
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("D14").Value = "'28.10"
$ws.Range("E14").Value = "  +3.63%  "
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "3.598.17"
$ws.Range("E16").Value = "  +1.71%  "
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "65.873.30"
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("E19").Value = "  -2.16%  "
$ws.Range("D20").Value = "'14.75"
$ws.Range("E20").Value = "  +2.90%  "
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").Value = "'398.46"
$ws.Range("E22").Value = "  +1.31%  "
$ws.Range("D23").Value = "'0.593"
$ws.Range("E23").Value = "  +3.30%  "
$ws.Range("D24").Value = "3.748.55"
$ws.Range("D25").Value = "'74.47"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +3.34%  "
$ws.Range("D28").Value = "'8.23"
$ws.Range("E28").Value = "  +6.26%  "
$ws.Range("E29").Value = "  +30.78%  "
$ws.Range("D30").Value = "'2.42"
$ws.Range("E30").Value = "  +5.56%  "
$ws.Range("D31").Value = "'8.61"
$ws.Range("E31").Value = "  +4.65%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("D33").Value = "3.609.21"
$ws.Range("E33").Value = "  +1.74%  "
$ws.Range("D34").Value = "'24.58"
$ws.Range("E34").Value = "  +3.24%  "
$ws.Range("E35").Value = "  +2.33%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").Value = "'5.41"
$ws.Range("E37").Value = "  +9.42%  "
$ws.Range("E38").Value = "  +4.06%  "
$ws.Range("D39").Value = "'7.11"
$ws.Range("E39").Value = "  +2.18%  "
$ws.Range("D40").Value = "'170.36"
$ws.Range("E40").Value = "  +0.96%  "
$ws.Range("D41").Value = "'0.0838"
$ws.Range("E41").Value = "  +3.49%  "
$ws.Range("D42").Value = "'0.844"
$ws.Range("E42").Value = "  +2.15%  "
$ws.Range("D43").Value = "'26.25"
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").Value = "'1.26"
$ws.Range("E44").Value = "  +6.66%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'43.35"
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("D46").Value = "'4.55"
$ws.Range("E46").Value = "  +2.67%  "
$ws.Range("E48").Value = "  +2.20%  "
$ws.Range("E49").Value = "  +4.16%  "
$ws.Range("D50").Value = "2.426.96"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "'320.20"
$ws.Range("E51").Value = "  +6.64%  "
